# ===== Step 0: get worksheet references (by current/original names) =====
$wb = $excel.ActiveWorkbook
$klWs = $wb.Worksheets.Item("kernel_launchers")        # sheet3
$oldXlaSummaryWs = $wb.Worksheets.Item("xla_summary")   # sheet7 (will become df_xla_perf)
$opGemmWs = $wb.Worksheets.Item("op_jax_gemm")          # sheet8 (will become xla_summary)
$opConvWs = $wb.Worksheets.Item("op_jax_conv")          # sheet9 (will become op_conv)
$opTeWs = $wb.Worksheets.Item("op_jax_te")              # sheet10 (will be deleted)

# ===== Step 1: snapshot old xla_summary content (A1:G4) into op_jax_gemm (sheet8) =====
# This preserves the original xla_summary table, which becomes the NEW xla_summary sheet (sheet8).
$oldXlaSummaryWs.Range("A1:G4").Copy($opGemmWs.Range("A1:G4"))

# ===== Step 2: add GPU_kernel_launch_latency column (M) to kernel_launchers (sheet3) =====
$klWs.Range("L1").Copy($klWs.Range("M1"))
$klWs.Cells.Item(1,13).Value = "GPU_kernel_launch_latency"
$klWs.Cells.Item(2,13).Value = 20.94200000003912
$klWs.Cells.Item(3,13).Value = 12.1589999999851
$klWs.Cells.Item(4,13).Value = 16.84600000013597
$klWs.Cells.Item(5,13).Value = 15.73399999993853
$klWs.Cells.Item(6,13).Value = 294.3659999999218
$klWs.Cells.Item(7,13).Value = 14.40200000000186
$klWs.Cells.Item(8,13).Value = 10.92700000014156
$klWs.Cells.Item(9,13).Value = 14.7019999998156
$klWs.Cells.Item(10,13).Value = 13.06099999998696
$klWs.Cells.Item(11,13).Value = 295.3469999998342
$klWs.Cells.Item(12,13).Value = 14.72200000006706
$klWs.Cells.Item(13,13).Value = 11.14699999988079
$klWs.Cells.Item(14,13).Value = 13.55099999997765
$klWs.Cells.Item(15,13).Value = 8.503000000026077
$klWs.Cells.Item(16,13).Value = 287.7750000001397
$klWs.Cells.Item(17,13).Value = 21.33199999993667
$klWs.Cells.Item(18,13).Value = 9.385000000009313
$klWs.Cells.Item(19,13).Value = 14.77299999981187
$klWs.Cells.Item(20,13).Value = 9.545000000158325
$klWs.Cells.Item(21,13).Value = 293.1029999998864
$klWs.Cells.Item(22,13).Value = 15.06200000015087
$klWs.Cells.Item(23,13).Value = 11.5570000000298
$klWs.Cells.Item(24,13).Value = 14.02099999994971
$klWs.Cells.Item(25,13).Value = 13.11999999987893
$klWs.Cells.Item(26,13).Value = 295.7669999999925

# ===== Step 3: rebuild xla_summary (sheet7) content as the new df_xla_perf table =====
$xlaWs = $oldXlaSummaryWs
$xlaWs.Cells.Clear()
$klWs.Range("A1:L1").Copy($xlaWs.Range("A1:L1"))
$xlaWs.Range("L1").Copy($xlaWs.Range("M1"))
$xlaWs.Cells.Item(1,13).Value = "GPU_kernel_launch_latency"
$xlaWs.Range("L1").Copy($xlaWs.Range("N1"))
$xlaWs.Cells.Item(1,14).Value = "total_input_bytes"

$srcRows = @(2,3,6,7,8,11,12,13,16,17,18,21,22,23,26)
$destRow = 2
foreach ($sr in $srcRows) {
    $klWs.Range("A" + $sr + ":L" + $sr).Copy($xlaWs.Range("A" + $destRow + ":L" + $destRow))
    $destRow++
}

$xlaWs.Cells.Item(2,13).Value = 20.94200000003912
$xlaWs.Cells.Item(2,14).Value = 655360
$xlaWs.Cells.Item(3,13).Value = 12.1589999999851
$xlaWs.Cells.Item(3,14).Value = 6389760
$xlaWs.Cells.Item(4,13).Value = 294.3659999999218
$xlaWs.Cells.Item(4,14).Value = 572026880
$xlaWs.Cells.Item(5,13).Value = 14.40200000000186
$xlaWs.Cells.Item(5,14).Value = 655360
$xlaWs.Cells.Item(6,13).Value = 10.92700000014156
$xlaWs.Cells.Item(6,14).Value = 6389760
$xlaWs.Cells.Item(7,13).Value = 295.3469999998342
$xlaWs.Cells.Item(7,14).Value = 572026880
$xlaWs.Cells.Item(8,13).Value = 14.72200000006706
$xlaWs.Cells.Item(8,14).Value = 655360
$xlaWs.Cells.Item(9,13).Value = 11.14699999988079
$xlaWs.Cells.Item(9,14).Value = 6389760
$xlaWs.Cells.Item(10,13).Value = 287.7750000001397
$xlaWs.Cells.Item(10,14).Value = 572026880
$xlaWs.Cells.Item(11,13).Value = 21.33199999993667
$xlaWs.Cells.Item(11,14).Value = 655360
$xlaWs.Cells.Item(12,13).Value = 9.385000000009313
$xlaWs.Cells.Item(12,14).Value = 6389760
$xlaWs.Cells.Item(13,13).Value = 293.1029999998864
$xlaWs.Cells.Item(13,14).Value = 572026880
$xlaWs.Cells.Item(14,13).Value = 15.06200000015087
$xlaWs.Cells.Item(14,14).Value = 655360
$xlaWs.Cells.Item(15,13).Value = 11.5570000000298
$xlaWs.Cells.Item(15,14).Value = 6389760
$xlaWs.Cells.Item(16,13).Value = 295.7669999999925
$xlaWs.Cells.Item(16,14).Value = 572026880

# ===== Step 4: fix op_jax_conv (sheet9) cell types: F2,G2,H2,F3,G3,H3 -> text =====
$opConvWs.Cells.Item(2,6).Value = "'False"
$opConvWs.Cells.Item(2,6).Style = "Normal"
$opConvWs.Cells.Item(2,7).Value = "'2"
$opConvWs.Cells.Item(2,7).Style = "Normal"
$opConvWs.Cells.Item(2,8).Value = "'False"
$opConvWs.Cells.Item(2,8).Style = "Normal"
$opConvWs.Cells.Item(3,6).Value = "'False"
$opConvWs.Cells.Item(3,6).Style = "Normal"
$opConvWs.Cells.Item(3,7).Value = "'2"
$opConvWs.Cells.Item(3,7).Style = "Normal"
$opConvWs.Cells.Item(3,8).Value = "'False"
$opConvWs.Cells.Item(3,8).Style = "Normal"

# ===== Step 5: delete op_jax_te (sheet10) =====
$opTeWs.Delete()

# ===== Step 6: rename sheets to final names =====
$oldXlaSummaryWs.Name = "df_xla_perf"   # sheet7: xla_summary -> df_xla_perf
$opGemmWs.Name = "xla_summary"          # sheet8: op_jax_gemm -> xla_summary
$opConvWs.Name = "op_conv"              # sheet9: op_jax_conv -> op_conv

